$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (columns A-T), reflecting the
# updated NATMI ligand-receptor computation results ("Natmi following
# Dr Hou advice"). Three sending clusters (ECs, FAPs, sCs) each paired
# with three target clusters (ECs, FAPs, sCs) for the Tfpi-Vldlr pair.
$data = @(
    @("ECs","Tfpi","Vldlr","ECs",3,1,4.288808333333333,12.866425,0.1993390083529519,0.1993390083529519,2,0.6666666666666666,0.3374003333333333,1.012201,0.01738364872808817,0.01738364872808818,1.447045361269444,13.023408251425,0.003465239299013151,0.003465239299013152),
    @("ECs","Tfpi","Vldlr","FAPs",3,1,4.288808333333333,12.866425,0.1993390083529519,0.1993390083529519,3,1,16.71131166666667,50.133935,0.8610055862391021,0.8610055862391023,71.67161273693056,645.044514632375,0.1716319997472546,0.1716319997472547),
    @("ECs","Tfpi","Vldlr","sCs",3,1,4.288808333333333,12.866425,0.1993390083529519,0.1993390083529519,3,1,2.360351,7.081053000000001,0.1216107650328097,0.1216107650328097,10.12309303839167,91.10783734552501,0.02424176930668412,0.02424176930668413),
    @("FAPs","Tfpi","Vldlr","ECs",3,1,11.10519266666667,33.315578,0.5161569185865866,0.5161569185865866,2,0.6666666666666666,0.3374003333333333,1.012201,0.01738364872808817,0.01738364872808818,3.746895707464222,33.722061367178,0.008972690561281629,0.008972690561281631),
    @("FAPs","Tfpi","Vldlr","FAPs",3,1,11.10519266666667,33.315578,0.5161569185865866,0.5161569185865866,3,1,16.71131166666667,50.133935,0.8610055862391021,0.8610055862391023,185.5823357710478,1670.24102193943,0.4444139902790125,0.4444139902790126),
    @("FAPs","Tfpi","Vldlr","sCs",3,1,11.10519266666667,33.315578,0.5161569185865866,0.5161569185865866,3,1,2.360351,7.081053000000001,0.1216107650328097,0.1216107650328097,26.21215261595934,235.909373543634,0.06277023774629245,0.06277023774629246),
    @("sCs","Tfpi","Vldlr","ECs",3,1,6.121147333333333,18.363442,0.2845040730604615,0.2845040730604615,2,0.6666666666666666,0.3374003333333333,1.012201,0.01738364872808817,0.01738364872808818,2.065277150649111,18.587494355842,0.004945718867793396,0.004945718867793398),
    @("sCs","Tfpi","Vldlr","FAPs",3,1,6.121147333333333,18.363442,0.2845040730604615,0.2845040730604615,3,1,16.71131166666667,50.133935,0.8610055862391021,0.8610055862391023,102.2924008449189,920.63160760427,0.244959596212835,0.244959596212835),
    @("sCs","Tfpi","Vldlr","sCs",3,1,6.121147333333333,18.363442,0.2845040730604615,0.2845040730604615,3,1,2.360351,7.081053000000001,0.1216107650328097,0.1216107650328097,14.44805622938067,130.032506064426,0.03459875797983309,0.0345987579798331)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
